# Update the "想去人数" (expected attendance) figures for two events that
# appear on both the "展览" sheet and the aggregated "全部类型" sheet.
#   F2: 4920 -> 4929
#   F4: 859  -> 864

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4929
    $ws.Range("F4").Value = 864
}
